# Update the cryptocurrency price/volume table with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.120.60"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "1.849.34"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "238.38"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.6933"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -6.17%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3057"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.93%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07579"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +4.71%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "23.45"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -4.94%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08111"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "1.862.35"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("E13").Value = "  -3.54%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.186"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.80%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "89.19"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").Value = "29.450.43"
$ws.Range("E16").Value = "  -1.62%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.798"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -5.07%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "241.94"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.96%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000007742"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "13.10"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -3.48%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.152.16"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.06%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.646"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -4.49%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.028"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1462"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -5.77%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "161.36"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("E28").Value = "  -3.05%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.939"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -4.12%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.392"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -7.06%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.433"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -3.90%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.497"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.69%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.049"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -4.82%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05241"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.87%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.188"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.76%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7112"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -5.16%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.663"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").Value = "  -5.01%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.9185"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +6.48%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.949"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -3.03%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.4280"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -5.37%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "70.01"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").Value = "1.045.34"
$ws.Range("E45").Value = "  -6.17%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "102.58"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "2.029.56"
$ws.Range("E48").Value = "  -0.39%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "7.231"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.83%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.747"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -6.25%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "9.268"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.55%  "
